# pcb progress: 60% complete
# Apply the "Pinouts" sheet edits: reshuffle the OUT_TLE / GPIO(EXTI) pin
# labels in column C/D (rows 10-17), rotate three MODULE-alias cells
# (rows 24,25,26,27,46), renumber a few STM_MOD labels (rows 34,38-40,50),
# and move the active selection to D35.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pinouts")

# New shared-string literals need to be introduced in the same order they
# appear in the target file's <sst> table (OUT_TLE6, OUT_TLE7, GPIO (EXTI2)).
$ws.Range("D34").Value = "OUT_TLE6"

# --- Column C / D block (rows 10-17): OUT_TLE* / GPIO(EXTI*) relabeling ---
$ws.Range("D10").Value = "OUT_TLE7"
$ws.Range("D11").Value = "OUT_TLE4"

$ws.Range("C12").Value = "GPIO (EXTI2)"
$ws.Range("D12").Value = "TLE_F"

$ws.Range("C13").Value = "GPIO"
$ws.Range("D13").Value = "OUT_TLE3"

$ws.Range("D14").Value = "OUT_TLE5"
$ws.Range("D15").Value = "OUT_TLE2"
$ws.Range("D17").Value = "OUT_TLE1"

# --- Rows 24/25/26/27: rotate the *MODULE2/3 + OUT5 combo out to row 27 ---
$ws.Range("O24").Value = ""
$ws.Range("P24").Value = "*MODULE2"
$ws.Range("P24").Font.Italic = $true
$ws.Range("T24").Value = ""

$ws.Range("O25").Value = ""
$ws.Range("P25").Value = "*MODULE3"
$ws.Range("P25").Font.Italic = $true
$ws.Range("T25").Value = ""

$ws.Range("P26").Font.Italic = $true

$ws.Range("O27").Value = "OUT5"
$ws.Range("P27").Value = "direct-5V, res for IGN5"
$ws.Range("P27").Font.Italic = $false
$ws.Range("T27").Value = 4

# --- Row 46: pick up the OUT5 / IGN6 pair that used to live on row 25 ---
$ws.Range("O46").Value = "OUT5"
$ws.Range("O46").Font.Italic = $true
$ws.Range("P46").Value = "direct-5V, res for IGN6"
$ws.Range("P46").Font.Italic = $false
$ws.Range("T46").Value = 5

# --- STM_MOD renumbering (rows 38-40, 50) ---
$ws.Range("D38").Value = "STM_MOD1"
$ws.Range("D39").Value = "STM_MOD2"
$ws.Range("D40").Value = "STM_MOD3"
$ws.Range("D50").Value = "STM_MOD4"

# --- Move the saved selection / scroll position ---
$ws.Activate()
$ws.Range("D35").Select()
